$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3 and 6 were highlighted (yellow fill) -- remove that highlight by
# pasting the plain formatting from row 2 (which already uses the default,
# unhighlighted cell style) onto them.
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6:F6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# "don mua" (order) quantities/totals update: quantity on row 5 and total on
# row 4 are zeroed out.
$ws.Range("F4").Value = 0
$ws.Range("E5").Value = 0

# Move the active cell selection to K6.
$ws.Range("K6").Select()
